# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets to reflect the latest scraped data.
#
# 展览  F2: 336 -> 335
# 展览  F3: 1391 -> 1395
# 展览  F4: 88  -> 89
# 全部类型 F2: 336 -> 335
# 全部类型 F3: 1391 -> 1395
# 全部类型 F4: 88  -> 89

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 335
    $ws.Range("F3").Value = 1395
    $ws.Range("F4").Value = 89
}
